$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new blank rows before row 67 (old row 67 shifts down to row 71)
$ws.Range("A67:T70").EntireRow.Insert()

# --- Row 65: update in place (Especial) ---
$ws.Cells.Item(65, 4).Value = 44448
$ws.Cells.Item(65, 13).Value = 125
$ws.Cells.Item(65, 14).Value = 3200
$ws.Cells.Item(65, 15).Value = 3200
$ws.Cells.Item(65, 16).Value = 3200
$ws.Cells.Item(65, 19).Value = 3200

# --- Row 66: update in place (Extra (doble especial)) ---
$ws.Cells.Item(66, 4).Value = 44448
$ws.Cells.Item(66, 13).Value = 100
$ws.Cells.Item(66, 14).Value = 3400
$ws.Cells.Item(66, 15).Value = 3400
$ws.Cells.Item(66, 16).Value = 3400
$ws.Cells.Item(66, 19).Value = 3400

# --- Row 67: new (Primera) ---
$ws.Cells.Item(67, 1).Value = 6
$ws.Cells.Item(67, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(67, 3).Value = "Metropolitana"
$ws.Cells.Item(67, 4).Value = 44448
$ws.Cells.Item(67, 5).Value = 13
$ws.Cells.Item(67, 6).Value = "Fruta"
$ws.Cells.Item(67, 7).Value = 100107
$ws.Cells.Item(67, 8).Value = "Otros"
$ws.Cells.Item(67, 9).Value = 100107002
$ws.Cells.Item(67, 10).Value = "Chirimoya"
$ws.Cells.Item(67, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(67, 12).Value = "Primera"
$ws.Cells.Item(67, 13).Value = 270
$ws.Cells.Item(67, 14).Value = 2800
$ws.Cells.Item(67, 15).Value = 3000
$ws.Cells.Item(67, 16).Value = 2900
$ws.Cells.Item(67, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(67, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(67, 19).Value = 2900
$ws.Cells.Item(67, 20).Value = 1

# --- Row 68: new (Segunda) ---
$ws.Cells.Item(68, 1).Value = 6
$ws.Cells.Item(68, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(68, 3).Value = "Metropolitana"
$ws.Cells.Item(68, 4).Value = 44448
$ws.Cells.Item(68, 5).Value = 13
$ws.Cells.Item(68, 6).Value = "Fruta"
$ws.Cells.Item(68, 7).Value = 100107
$ws.Cells.Item(68, 8).Value = "Otros"
$ws.Cells.Item(68, 9).Value = 100107002
$ws.Cells.Item(68, 10).Value = "Chirimoya"
$ws.Cells.Item(68, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(68, 12).Value = "Segunda"
$ws.Cells.Item(68, 13).Value = 230
$ws.Cells.Item(68, 14).Value = 2500
$ws.Cells.Item(68, 15).Value = 2500
$ws.Cells.Item(68, 16).Value = 2500
$ws.Cells.Item(68, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(68, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(68, 19).Value = 2500
$ws.Cells.Item(68, 20).Value = 1

# --- Row 69: new (Especial, matches old row65 data, date 44167) ---
$ws.Cells.Item(69, 1).Value = 6
$ws.Cells.Item(69, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(69, 3).Value = "Metropolitana"
$ws.Cells.Item(69, 4).Value = 44167
$ws.Cells.Item(69, 5).Value = 13
$ws.Cells.Item(69, 6).Value = "Fruta"
$ws.Cells.Item(69, 7).Value = 100107
$ws.Cells.Item(69, 8).Value = "Otros"
$ws.Cells.Item(69, 9).Value = 100107002
$ws.Cells.Item(69, 10).Value = "Chirimoya"
$ws.Cells.Item(69, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(69, 12).Value = "Especial"
$ws.Cells.Item(69, 13).Value = 450
$ws.Cells.Item(69, 14).Value = 1700
$ws.Cells.Item(69, 15).Value = 1700
$ws.Cells.Item(69, 16).Value = 1700
$ws.Cells.Item(69, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(69, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(69, 19).Value = 1700
$ws.Cells.Item(69, 20).Value = 1

# --- Row 70: new (Extra (doble especial), matches old row66 data, date 44167) ---
$ws.Cells.Item(70, 1).Value = 6
$ws.Cells.Item(70, 2).Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Cells.Item(70, 3).Value = "Metropolitana"
$ws.Cells.Item(70, 4).Value = 44167
$ws.Cells.Item(70, 5).Value = 13
$ws.Cells.Item(70, 6).Value = "Fruta"
$ws.Cells.Item(70, 7).Value = 100107
$ws.Cells.Item(70, 8).Value = "Otros"
$ws.Cells.Item(70, 9).Value = 100107002
$ws.Cells.Item(70, 10).Value = "Chirimoya"
$ws.Cells.Item(70, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(70, 12).Value = "Extra (doble especial)"
$ws.Cells.Item(70, 13).Value = 300
$ws.Cells.Item(70, 14).Value = 2200
$ws.Cells.Item(70, 15).Value = 2200
$ws.Cells.Item(70, 16).Value = 2200
$ws.Cells.Item(70, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(70, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(70, 19).Value = 2200
$ws.Cells.Item(70, 20).Value = 1

# Row 71 already holds the shifted former row-67 data (Primera, date 44167) unchanged.
